# Remove the "ECs" target-cluster rows (old rows 8, 5, 2) and renumber
# the remaining rows, then refresh every data cell with the recalculated
# TPM-derived values (ECs removed as both a sending and a receiving
# cluster changes the specificity normalisation for every remaining row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so row numbers of not-yet-deleted rows don't shift.
$ws.Range("A8:T8").EntireRow.Delete()
$ws.Range("A5:T5").EntireRow.Delete()
$ws.Range("A2:T2").EntireRow.Delete()

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bsg"
$ws.Range("C2").Value = "Slc16a7"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.84781966666667
$ws.Range("H2").Value = 101.543459
$ws.Range("I2").Value = 0.2402182618707165
$ws.Range("J2").Value = 0.2402182618707166
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.286000666666667
$ws.Range("N2").Value = 6.858002
$ws.Range("O2").Value = 0.2613437960247708
$ws.Range("P2").Value = 0.2613437960247707
$ws.Range("Q2").Value = 77.37613832321311
$ws.Range("R2").Value = 696.385244908918
$ws.Range("S2").Value = 0.06277955243176551
$ws.Range("T2").Value = 0.0627795524317655
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bsg"
$ws.Range("C3").Value = "Slc16a7"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.84781966666667
$ws.Range("H3").Value = 101.543459
$ws.Range("I3").Value = 0.2402182618707165
$ws.Range("J3").Value = 0.2402182618707166
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.461100666666667
$ws.Range("N3").Value = 19.383302
$ws.Range("O3").Value = 0.7386562039752294
$ws.Range("P3").Value = 0.7386562039752294
$ws.Range("Q3").Value = 218.6941702135131
$ws.Range("R3").Value = 1968.247531921618
$ws.Range("S3").Value = 0.1774387094389511
$ws.Range("T3").Value = 0.1774387094389511
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bsg"
$ws.Range("C4").Value = "Slc16a7"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 74.45592499999999
$ws.Range("H4").Value = 223.367775
$ws.Range("I4").Value = 0.5284143281787288
$ws.Range("J4").Value = 0.5284143281787288
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.286000666666667
$ws.Range("N4").Value = 6.858002
$ws.Range("O4").Value = 0.2613437960247708
$ws.Range("P4").Value = 0.2613437960247707
$ws.Range("Q4").Value = 170.2062941872833
$ws.Range("R4").Value = 1531.85664768555
$ws.Range("S4").Value = 0.138097806400108
$ws.Range("T4").Value = 0.1380978064001079
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bsg"
$ws.Range("C5").Value = "Slc16a7"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 74.45592499999999
$ws.Range("H5").Value = 223.367775
$ws.Range("I5").Value = 0.5284143281787288
$ws.Range("J5").Value = 0.5284143281787288
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.461100666666667
$ws.Range("N5").Value = 19.383302
$ws.Range("O5").Value = 0.7386562039752294
$ws.Range("P5").Value = 0.7386562039752294
$ws.Range("Q5").Value = 481.0672266547833
$ws.Range("R5").Value = 4329.60503989305
$ws.Range("S5").Value = 0.3903165217786209
$ws.Range("T5").Value = 0.3903165217786209
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bsg"
$ws.Range("C6").Value = "Slc16a7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 32.60069533333333
$ws.Range("H6").Value = 97.802086
$ws.Range("I6").Value = 0.2313674099505547
$ws.Range("J6").Value = 0.2313674099505547
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.286000666666667
$ws.Range("N6").Value = 6.858002
$ws.Range("O6").Value = 0.2613437960247708
$ws.Range("P6").Value = 0.2613437960247707
$ws.Range("Q6").Value = 74.5252112657969
$ws.Range("R6").Value = 670.7269013921721
$ws.Range("S6").Value = 0.06046643719289728
$ws.Range("T6").Value = 0.06046643719289727
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bsg"
$ws.Range("C7").Value = "Slc16a7"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 32.60069533333333
$ws.Range("H7").Value = 97.802086
$ws.Range("I7").Value = 0.2313674099505547
$ws.Range("J7").Value = 0.2313674099505547
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.461100666666667
$ws.Range("N7").Value = 19.383302
$ws.Range("O7").Value = 0.7386562039752294
$ws.Range("P7").Value = 0.7386562039752294
$ws.Range("Q7").Value = 210.6363743519969
$ws.Range("R7").Value = 1895.727369167972
$ws.Range("S7").Value = 0.1709009727576574
$ws.Range("T7").Value = 0.1709009727576574
